$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.828084232555808
$ws.Range("C2").Value = 4.631832391164917
$ws.Range("D2").Value = 14.8832573780194
$ws.Range("E2").Value = 16.26228367436907
$ws.Range("G2").Value = 38.17157358099266
$ws.Range("H2").Value = 16.7364713804265
$ws.Range("J2").Value = 9.465869602604881
$ws.Range("K2").Value = 9.179785395928274
$ws.Range("M2").Value = 16.06445695289349
$ws.Range("O2").Value = 26.71244451377827
$ws.Range("B3").Value = 9.57701626589658
$ws.Range("C3").Value = 4.501380056928865
$ws.Range("D3").Value = 14.85805667283165
$ws.Range("E3").Value = 16.2633032609611
$ws.Range("G3").Value = 38.2671052148101
$ws.Range("H3").Value = 16.786805880836
$ws.Range("J3").Value = 9.48440005652667
$ws.Range("K3").Value = 9.012101957923397
$ws.Range("M3").Value = 16.00503280277584
$ws.Range("O3").Value = 26.79479281232389
$ws.Range("B4").Value = 9.421110697171155
$ws.Range("C4").Value = 4.418745397307377
$ws.Range("D4").Value = 14.84556422182238
$ws.Range("E4").Value = 16.26691877207267
$ws.Range("G4").Value = 38.3356770413121
$ws.Range("H4").Value = 16.82017123135026
$ws.Range("J4").Value = 9.49672383132679
$ws.Range("K4").Value = 8.908844442835157
$ws.Range("M4").Value = 15.97097382491599
$ws.Range("O4").Value = 26.85039077848597
$ws.Range("B5").Value = 9.357234676550924
$ws.Range("C5").Value = 4.384467034656044
$ws.Range("D5").Value = 14.8412271872126
$ws.Range("E5").Value = 16.26914497325771
$ws.Range("G5").Value = 38.36610604380594
$ws.Range("H5").Value = 16.83438648085941
$ws.Range("J5").Value = 9.501984082269509
$ws.Range("K5").Value = 8.866744886655404
$ws.Range("M5").Value = 15.9577149984069
$ws.Range("O5").Value = 26.87431197372583
$ws.Range("B6").Value = 9.346610293214017
$ws.Range("C6").Value = 4.378739693570344
$ws.Range("D6").Value = 14.84055266032366
$ws.Range("E6").Value = 16.26956012904296
$ws.Range("G6").Value = 38.37130863016466
$ws.Range("H6").Value = 16.83678427614432
$ws.Range("J6").Value = 9.502871940141615
$ws.Range("K6").Value = 8.859754598810454
$ws.Range("M6").Value = 15.95555113630626
$ws.Range("O6").Value = 26.878360395708
$ws.Range("B7").Value = 9.420250501596838
$ws.Range("C7").Value = 4.418285506601811
$ws.Range("D7").Value = 14.84550267411315
$ws.Range("E7").Value = 16.26694574589104
$ws.Range("G7").Value = 38.33607736487404
$ws.Range("H7").Value = 16.82036043840991
$ws.Range("J7").Value = 9.496793807955886
$ws.Range("K7").Value = 8.908276687620431
$ws.Range("M7").Value = 15.97079248666909
$ws.Range("O7").Value = 26.85070827021404
$ws.Range("B8").Value = 9.741935385276221
$ws.Range("C8").Value = 4.58739824828707
$ws.Range("D8").Value = 14.8739521696467
$ws.Range("E8").Value = 16.26201552626413
$ws.Range("G8").Value = 38.20245053227692
$ws.Range("H8").Value = 16.75331623634333
$ws.Range("J8").Value = 9.472062761238714
$ws.Range("K8").Value = 9.12206227406878
$ws.Range("M8").Value = 16.04346983766725
$ws.Range("O8").Value = 26.73979202022405
$ws.Range("B9").Value = 10.35489539329951
$ws.Range("C9").Value = 4.897559376525868
$ws.Range("D9").Value = 14.95318536312966
$ws.Range("E9").Value = 16.27600641920017
$ws.Range("G9").Value = 38.01940177376478
$ws.Range("H9").Value = 16.64135663601938
$ws.Range("J9").Value = 9.431056576168222
$ws.Range("K9").Value = 9.536613684032323
$ws.Range("M9").Value = 16.20478552183949
$ws.Range("O9").Value = 26.56231692175768
$ws.Range("B10").Value = 10.78939629265125
$ws.Range("C10").Value = 5.11073761526495
$ws.Range("D10").Value = 15.02537525006397
$ws.Range("E10").Value = 16.30060907667511
$ws.Range("G10").Value = 37.93347231809921
$ws.Range("H10").Value = 16.57099150159504
$ws.Range("J10").Value = 9.40547616107604
$ws.Range("K10").Value = 9.835354722312212
$ws.Range("M10").Value = 16.33409982536489
$ws.Range("O10").Value = 26.45642926270238
$ws.Range("B11").Value = 10.98269227467667
$ws.Range("C11").Value = 5.204237648052987
$ws.Range("D11").Value = 15.06117224995751
$ws.Range("E11").Value = 16.31488347588565
$ws.Range("G11").Value = 37.90499482944883
$ws.Range("H11").Value = 16.54156182659083
$ws.Range("J11").Value = 9.394821995250428
$ws.Range("K11").Value = 9.969403332459621
$ws.Range("M11").Value = 16.39511988930943
$ws.Range("O11").Value = 26.41359873924459
$ws.Range("B12").Value = 11.05519202204984
$ws.Range("C12").Value = 5.239122922194621
$ws.Range("D12").Value = 15.07514530538662
$ws.Range("E12").Value = 16.32072906765034
$ws.Range("G12").Value = 37.89574159588679
$ws.Range("H12").Value = 16.53078851919293
$ws.Range("J12").Value = 9.390928483067766
$ws.Range("K12").Value = 10.01985187825481
$ws.Range("M12").Value = 16.41852846358958
$ws.Range("O12").Value = 26.39814910724809
$ws.Range("B13").Value = 11.03960996833011
$ws.Range("C13").Value = 5.231633248887355
$ws.Range("D13").Value = 15.07211751858891
$ws.Range("E13").Value = 16.31945059022981
$ws.Range("G13").Value = 37.89766630941722
$ws.Range("H13").Value = 16.53309223855893
$ws.Range("J13").Value = 9.391760754717431
$ws.Range("K13").Value = 10.00900150737789
$ws.Range("M13").Value = 16.4134738211456
$ws.Range("O13").Value = 26.40144221822106
$ws.Range("B14").Value = 10.98867122708968
$ws.Range("C14").Value = 5.207118239609828
$ws.Range("D14").Value = 15.06231350058403
$ws.Range("E14").Value = 16.31535559312983
$ws.Range("G14").Value = 37.90420286240891
$ws.Range("H14").Value = 16.54066806269232
$ws.Range("J14").Value = 9.394498849705611
$ws.Range("K14").Value = 9.973560291117591
$ws.Range("M14").Value = 16.39703976369206
$ws.Range("O14").Value = 26.41231226015646
$ws.Range("B15").Value = 10.95737699459516
$ws.Range("C15").Value = 5.192033623345147
$ws.Range("D15").Value = 15.05636238849084
$ws.Range("E15").Value = 16.31290451706925
$ws.Range("G15").Value = 37.90840613591789
$ws.Range("H15").Value = 16.54535680125184
$ws.Range("J15").Value = 9.396194364204133
$ws.Range("K15").Value = 9.951809450724996
$ws.Range("M15").Value = 16.38701228498202
$ws.Range("O15").Value = 26.41907071710398
$ws.Range("B16").Value = 10.77667006076955
$ws.Range("C16").Value = 5.10455534982338
$ws.Range("D16").Value = 15.02309469183456
$ws.Range("E16").Value = 16.29973798452823
$ws.Range("G16").Value = 37.93554734256109
$ws.Range("H16").Value = 16.57296672062359
$ws.Range("J16").Value = 9.406192181207226
$ws.Range("K16").Value = 9.826553049937067
$ws.Range("M16").Value = 16.33015499466519
$ws.Range("O16").Value = 26.45933590729958
$ws.Range("B17").Value = 10.66464348828775
$ws.Range("C17").Value = 5.049984113218504
$ws.Range("D17").Value = 15.00343815974329
$ws.Range("E17").Value = 16.29244817609629
$ws.Range("G17").Value = 37.95491914951468
$ws.Range("H17").Value = 16.59056532852622
$ws.Range("J17").Value = 9.41257695082315
$ws.Range("K17").Value = 9.749204399501622
$ws.Range("M17").Value = 16.2958271341541
$ws.Range("O17").Value = 26.48540581563761
$ws.Range("B18").Value = 10.599802132017
$ws.Range("C18").Value = 5.018270147517337
$ws.Range("D18").Value = 14.9924111128936
$ws.Range("E18").Value = 16.28854561989221
$ws.Range("G18").Value = 37.96705995724881
$ws.Range("H18").Value = 16.60093039708184
$ws.Range("J18").Value = 9.416341795912276
$ws.Range("K18").Value = 9.704543989159008
$ws.Range("M18").Value = 16.27629006238955
$ws.Range("O18").Value = 26.50090290690388
$ws.Range("B19").Value = 10.57778036033851
$ws.Range("C19").Value = 5.007477034589006
$ws.Range("D19").Value = 14.98872566833512
$ws.Range("E19").Value = 16.28727423702821
$ws.Range("G19").Value = 37.9713419940185
$ws.Range("H19").Value = 16.60448153172074
$ws.Range("J19").Value = 9.417632402664736
$ws.Range("K19").Value = 9.689394748594919
$ws.Range("M19").Value = 16.26971118520402
$ws.Range("O19").Value = 26.50623618438334
$ws.Range("B20").Value = 10.67661152085009
$ws.Range("C20").Value = 5.055827219547464
$ws.Range("D20").Value = 15.00550182364355
$ws.Range("E20").Value = 16.29319415980189
$ws.Range("G20").Value = 37.95275359470667
$ws.Range("H20").Value = 16.58866679549361
$ws.Range("J20").Value = 9.411887710477519
$ws.Range("K20").Value = 9.757456392626471
$ws.Range("M20").Value = 16.2994600273924
$ws.Range("O20").Value = 26.48257862328187
$ws.Range("B21").Value = 11.00365262285501
$ws.Range("C21").Value = 5.214333195133309
$ws.Range("D21").Value = 15.06518191188567
$ws.Range("E21").Value = 16.31654647278782
$ws.Range("G21").Value = 37.90224134706403
$ws.Range("H21").Value = 16.53843278790833
$ws.Range("J21").Value = 9.393690780765249
$ws.Range("K21").Value = 9.983979082956433
$ws.Range("M21").Value = 16.40185877266012
$ws.Range("O21").Value = 26.40909857199209
$ws.Range("B22").Value = 11.21330044813841
$ws.Range("C22").Value = 5.314879955152458
$ws.Range("D22").Value = 15.10661587564984
$ws.Range("E22").Value = 16.33437262625364
$ws.Range("G22").Value = 37.87815152253471
$ws.Range("H22").Value = 16.50776489701824
$ws.Range("J22").Value = 9.382619691476865
$ws.Range("K22").Value = 10.13018195274255
$ws.Range("M22").Value = 16.47053346578414
$ws.Range("O22").Value = 26.36556024434429
$ws.Range("B23").Value = 11.10180356234872
$ws.Range("C23").Value = 5.261501447989598
$ws.Range("D23").Value = 15.08428218930584
$ws.Range("E23").Value = 16.32462496846977
$ws.Range("G23").Value = 37.89019099573466
$ws.Range("H23").Value = 16.52393498514058
$ws.Range("J23").Value = 9.388453455523511
$ws.Range("K23").Value = 10.05233410823497
$ws.Range("M23").Value = 16.43372497155046
$ws.Range("O23").Value = 26.38838653697459
$ws.Range("B24").Value = 10.67120212531216
$ws.Range("C24").Value = 5.053186609307959
$ws.Range("D24").Value = 15.00456798790021
$ws.Range("E24").Value = 16.29285600184713
$ws.Range("G24").Value = 37.95372951590655
$ws.Range("H24").Value = 16.58952435165146
$ws.Range("J24").Value = 9.412199022726542
$ws.Range("K24").Value = 9.753726259873318
$ws.Range("M24").Value = 16.29781697708423
$ws.Range("O24").Value = 26.48385521107129
$ws.Range("B25").Value = 10.19151956505997
$ws.Range("C25").Value = 4.816127214969923
$ws.Range("D25").Value = 14.92927100499099
$ws.Range("E25").Value = 16.26969755680049
$ws.Range("G25").Value = 38.06042032149776
$ws.Range("H25").Value = 16.66955605856261
$ws.Range("J25").Value = 9.441349909839463
$ws.Range("K25").Value = 9.425282384814674
$ws.Range("M25").Value = 16.1591992800226
$ws.Range("O25").Value = 26.6060324741689
